$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contains one price record per row (rows 2-207). A new weekly
# record is inserted as row 101, pushing the existing rows 101-207 down to
# 102-208 (dimension grows from T207 to T208).
$ws.Rows.Item(101).Insert()

# Populate the newly inserted row with the new record's data. Columns
# A,B,C,E,F,G,H,I,J,R share the same "Macroferia Regional de Talca / Piña
# / Ecuador" values used throughout the sheet.
$ws.Cells.Item(101, 1).Value = 5
$ws.Cells.Item(101, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(101, 3).Value = "Maule"
$ws.Cells.Item(101, 4).Value = 44601
$ws.Cells.Item(101, 5).Value = 7
$ws.Cells.Item(101, 6).Value = "Fruta"
$ws.Cells.Item(101, 7).Value = 100108
$ws.Cells.Item(101, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(101, 9).Value = 100108005
$ws.Cells.Item(101, 10).Value = "Piña"
$ws.Cells.Item(101, 11).Value = "Caramelo"
$ws.Cells.Item(101, 12).Value = "Tercera"
$ws.Cells.Item(101, 13).Value = 260
$ws.Cells.Item(101, 14).Value = 15000
$ws.Cells.Item(101, 15).Value = 15000
$ws.Cells.Item(101, 16).Value = 15000
$ws.Cells.Item(101, 17).Value = "$/caja 16 unidades"
$ws.Cells.Item(101, 18).Value = "Ecuador"
$ws.Cells.Item(101, 19).Value = 938
$ws.Cells.Item(101, 20).Value = 16
